# Added ErrorMessages constants file
#
# Appends four new student records to the "student" sheet (sheet1):
#   id=3, bb,  dd, ff, 12246578, bb@bb.in,     bb123
#   id=4, cc,  dd, FG, 1346798,  cc@gmail.com, ccc
#   id=4, cc,  dd, FG, 1346798,  cc@gmail.com, ccc
#   id=6, tt,  dd, FG, 1346798,  tt@gmail.com, ttt
#
# Columns (row 1 header): id | firstname | lastname | qualification | mobile | email | password

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The new rows live directly below the existing data (rows 2-3), starting at row 4.
$startRow = 4

$newRecords = @(
    @("3", "bb", "dd", "ff", "12246578", "bb@bb.in", "bb123"),
    @("4", "cc", "dd", "FG", "1346798", "cc@gmail.com", "ccc"),
    @("4", "cc", "dd", "FG", "1346798", "cc@gmail.com", "ccc"),
    @("6", "tt", "dd", "FG", "1346798", "tt@gmail.com", "ttt")
)

$endRow = $startRow + $newRecords.Count - 1

# Pre-format the target block as Text so that numeric-looking values (ids, mobile
# numbers) are stored as literal strings instead of being auto-coerced to numbers,
# matching how the original data rows (2-3) are stored.
$targetRange = $ws.Range("A$($startRow):G$($endRow)")
$targetRange.NumberFormat = "@"

for ($i = 0; $i -lt $newRecords.Count; $i++) {
    $row = $startRow + $i
    $record = $newRecords[$i]
    for ($col = 1; $col -le $record.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $record[$col - 1]
    }
}

# Restore the default (Normal) style on the newly written cells so they don't carry
# an explicit style index, matching the look of the rest of the sheet.
$targetRange.Style = "Normal"
